# Weekly update: insert a new price observation for "Espinaca" at
# Vega Modelo de Temuco, pushing the existing rows 137:172 down to 138:173.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137 (shifts 137..172 -> 138..173, keeps formatting)
$ws.Rows(137).Insert()

# Populate the newly inserted row with the new weekly observation
$ws.Cells.Item(137, 1).Value  = 10
$ws.Cells.Item(137, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(137, 3).Value  = "La Araucanía"
$ws.Cells.Item(137, 4).Value  = 44785
$ws.Cells.Item(137, 5).Value  = 9
$ws.Cells.Item(137, 6).Value  = 100112012
$ws.Cells.Item(137, 7).Value  = "Espinaca"
$ws.Cells.Item(137, 8).Value  = "Sin especificar"
$ws.Cells.Item(137, 9).Value  = "Primera"
$ws.Cells.Item(137, 10).Value = 20
$ws.Cells.Item(137, 11).Value = 13000
$ws.Cells.Item(137, 12).Value = 13000
$ws.Cells.Item(137, 13).Value = 13000
$ws.Cells.Item(137, 14).Value = "$/docena de atados"
$ws.Cells.Item(137, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(137, 16).Value = 4333
$ws.Cells.Item(137, 17).Value = 3
$ws.Cells.Item(137, 18).Value = "Hortaliza"
